$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Passwords in D4:D12 all become 111111
$ws.Range("D4").Value = 111111
$ws.Range("D5").Value = 111111
$ws.Range("D6").Value = 111111
$ws.Range("D7").Value = 111111
$ws.Range("D8").Value = 111111
$ws.Range("D9").Value = 111111
$ws.Range("D10").Value = 111111
$ws.Range("D11").Value = 111111
$ws.Range("D12").Value = 111111

# 2. C5 username changes from "adadi" to "asadi"
$ws.Range("C5").Value = "asadi"

# 3. New rows 13-16: Code column (A) first
$ws.Range("A13").Value = "E101"
$ws.Range("A14").Value = "E102"
$ws.Range("A15").Value = "E103"
$ws.Range("A16").Value = "E104"

# 4. New rows 13-16: Name column (B)
$ws.Range("B13").Value = "خانم مسلمی"
$ws.Range("B14").Value = "خانم موسوی"
$ws.Range("B15").Value = "خانم شوشتری"
$ws.Range("B16").Value = "خانم فروغ نيا"

# 5. New row 17: Code then Name
$ws.Range("A17").Value = "A1001"
$ws.Range("B17").Value = "آقای دهنوی "

# 6. New rows 13-17: Type column (E)
$ws.Range("E13").Value = "Edari"
$ws.Range("E14").Value = "Edari"
$ws.Range("E15").Value = "Edari"
$ws.Range("E16").Value = "Edari"
$ws.Range("E17").Value = "Anbar"

# 7. New rows 13-17: Username column (C)
$ws.Range("C13").Value = "moslemi"
$ws.Range("C14").Value = "Kmosavi"
$ws.Range("C15").Value = "Eshoshtari"
$ws.Range("C16").Value = "Foroghniya"
$ws.Range("C17").Value = "dehnavi"

# 8. New rows 13-17: Password column (D)
$ws.Range("D13").Value = 111111
$ws.Range("D14").Value = 111111
$ws.Range("D15").Value = 111111
$ws.Range("D16").Value = 111111
$ws.Range("D17").Value = 111111

# 9. Apply the text/code number format (matching existing A-column style) to the
#    newly added code cells so they keep their leading characters / formatting.
$ws.Range("A13:A17").NumberFormat = "@"

# 10. Column C needs a custom width, same as the author resized it to fit the
#     longer usernames that were just added.
$ws.Columns.Item(3).ColumnWidth = 15.035714285714286

# 11. Move / collapse the selection to where the author left off after entering
#     the new rows.
$ws.Range("C18").Select()
